# "run on more subjects"
# The "In Lab" sheet had an extra (empty) column I that separated the
# timing table (A:G) from the supplementary results table that used to
# live in columns J:M. Now that more subjects are being run, that blank
# spacer column is removed so the results table shifts left to I:L,
# keeping everything compact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# Delete the whole, empty column I - this shifts J->I, K->J, L->K, M->L
$ws.Columns.Item(9).Delete()

# Restore the view state reported after the edit
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("H22").Select()
